$d = $word.ActiveDocument
$end = $d.Content.End
$rng = $d.Range($end, $end)

$bodyXml = @'
<w:p/><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>Table Sprint Planning</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="9634" w:type="dxa"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="562"/><w:gridCol w:w="9072"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9634" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="000000" w:themeFill="text1"/></w:tcPr><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Sprint Planning</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9634" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Sprint 1</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>NO</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Sprint Goals</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">Pendengar musik </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>dapat</w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> melakukan pendaftaran akun</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve">Pendengar musik </w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>dapat</w:t></w:r><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> melakukan login setelah melakukan register</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Pendengar musik dapat melihat beberapa rekomendasi musik</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Pendengar musik dapat melakukan pencarian musik</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Pendengar musik dapat melihat topchart yang berisi musik-musik yang popular</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Pendengar musik ingin jika musik di klik maka musik akan berpindah ke halaman yang lebih detail</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>7</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Pendengar musik ingin bisa melakukan play musik agar musik bisa di dengarkan</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>8</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Pendengar musik ingin bisa melakukan pause musik agar musik dapat berhenti</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9634" w:type="dxa"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/></w:tcPr><w:p><w:pPr><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="id-ID"/></w:rPr><w:t>Sprint Backlog</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>ID</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>Backlog Items</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>111</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(3) Sebagai pendengar musik saya ingin bisa mendaftarkan akun agar akun saya tidak hilang</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>222</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(3) Sebagai pendengar musik saya ingin bisa melakukan login akun setelah melakukan register</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>333</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(5) Sebagai pendengar musik saya ingin dibagian homepage terdapat fitur rekomendasi musik</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>444</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(13) Sebagai pendengar musik saya ingin dapat melakukan pencarian musik</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>555</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(5) Sebagai pendengar musik saya ingin di bagian hompage terdapat topchart yang berisi musik-musik yang sedang popular</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>666</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(8) Sebagai pendengar musik saya ingin jika musik di klik maka musik akan berpindah kehalaman musik yang lebih detail</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>777</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>13) Sebagai pendengar musik saya ingin bisa melakukan play musik agar musik dapat di dengarkan</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="562" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>888</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="9072" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="id-ID"/></w:rPr><w:t>(13) Sebagai pendengar musik saya ingin bisa melakukan pause musik agar musik dapat berhenti berputar</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:rPr><w:lang w:val="id-ID"/></w:rPr></w:pPr></w:p>
'@

$fullXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + "`n" + '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($fullXml)

Write-Host "Done. Paragraphs:" $d.Paragraphs.Count "Tables:" $d.Tables.Count
